$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031099476996261
$ws.Range("D2").Value = 1.033059460869538
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.039549971794471
$ws.Range("I2").Value = 1.032195244410111
$ws.Range("J2").Value = 1.036237070836587
$ws.Range("K2").Value = 1.03586290053958
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.042334841754393
$ws.Range("N2").Value = 1.037708646335244
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03266120905759
$ws.Range("D3").Value = 1.034189247763477
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.041198672110288
$ws.Range("I3").Value = 1.032537926270266
$ws.Range("J3").Value = 1.037437417667727
$ws.Range("K3").Value = 1.03680073182416
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.043791575183113
$ws.Range("N3").Value = 1.03891069779655
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033669303318899
$ws.Range("D4").Value = 1.034918109417741
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.042263340537683
$ws.Range("I4").Value = 1.032757334002351
$ws.Range("J4").Value = 1.038211353811081
$ws.Range("K4").Value = 1.037404813188469
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.044731568400971
$ws.Range("N4").Value = 1.039685733017999
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034092530392379
$ws.Range("D5").Value = 1.035224006349702
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.042710423556299
$ws.Range("I5").Value = 1.032849017922949
$ws.Range("J5").Value = 1.038536062160324
$ws.Range("K5").Value = 1.037658115128601
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.045126126846363
$ws.Range("N5").Value = 1.040010902490343
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034163558635067
$ws.Range("D6").Value = 1.035275337699553
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.042785461559744
$ws.Range("I6").Value = 1.032864379586713
$ws.Range("J6").Value = 1.038590543988537
$ws.Range("K6").Value = 1.037700607482237
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.045192339196916
$ws.Range("N6").Value = 1.040065461689001
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033674960746918
$ws.Range("D7").Value = 1.034922198847507
$ws.Range("E7").Value = 0.9943035907978915
$ws.Range("F7").Value = 1.042269316445937
$ws.Range("I7").Value = 1.032758561264584
$ws.Range("J7").Value = 1.03821569513967
$ws.Range("K7").Value = 1.037408200381651
$ws.Range("L7").Value = 0.9968970624459041
$ws.Range("M7").Value = 1.044736842911336
$ws.Range("N7").Value = 1.039690080511773
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031627785747098
$ws.Range("D8").Value = 1.033441734397106
$ws.Range("E8").Value = 0.9929600610674294
$ws.Range("F8").Value = 1.040107610003861
$ws.Range("I8").Value = 1.032311539556074
$ws.Range("J8").Value = 1.036643311976479
$ws.Range("K8").Value = 1.036180420014114
$ws.Range("L8").Value = 0.9958175282591053
$ws.Range("M8").Value = 1.042827698432624
$ws.Range("N8").Value = 1.038115464384135
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028001137084611
$ws.Range("D9").Value = 1.03081591638679
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.036281422448705
$ws.Range("I9").Value = 1.031505855984688
$ws.Range("J9").Value = 1.033850990298999
$ws.Range("K9").Value = 1.033995498571133
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.039443082463413
$ws.Range("N9").Value = 1.035319177289589
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025569723903455
$ws.Range("D10").Value = 1.029053472604366
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.033718493103844
$ws.Range("I10").Value = 1.030956469968395
$ws.Range("J10").Value = 1.031974400111391
$ws.Range("K10").Value = 1.032524075292664
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.037172263087407
$ws.Range("N10").Value = 1.033439922128668
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024513502550762
$ws.Range("D11").Value = 1.028287395052013
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.032605676725074
$ws.Range("I11").Value = 1.030715630640543
$ws.Range("J11").Value = 1.031158128570329
$ws.Range("K11").Value = 1.031883327998885
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.03618540938221
$ws.Range("N11").Value = 1.032622491388389
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024120650081925
$ws.Range("D12").Value = 1.02800239134057
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.0321918548374
$ws.Range("I12").Value = 1.030625725204026
$ws.Range("J12").Value = 1.030854363355786
$ws.Range("K12").Value = 1.031644775309793
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.035818298327426
$ws.Range("N12").Value = 1.032318294792399
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024204942238097
$ws.Range("D13").Value = 1.028063546033208
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.032280642636774
$ws.Range("I13").Value = 1.030645030509857
$ws.Range("J13").Value = 1.030919547789059
$ws.Range("K13").Value = 1.031695970718196
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.035897069973061
$ws.Range("N13").Value = 1.032383571795045
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024481040028158
$ws.Range("D14").Value = 1.028263845753244
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.032571479795479
$ws.Range("I14").Value = 1.030708208168187
$ws.Range("J14").Value = 1.031133030841952
$ws.Range("K14").Value = 1.031863620438501
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.036155075154557
$ws.Range("N14").Value = 1.03259735801836
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024651083037384
$ws.Range("D15").Value = 1.028387197347261
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.03275061115884
$ws.Range("I15").Value = 1.030747074671647
$ws.Range("J15").Value = 1.031264489533606
$ws.Range("K15").Value = 1.031966841759899
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.036313967445877
$ws.Range("N15").Value = 1.032729003396434
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025639748980564
$ws.Range("D16").Value = 1.029104252237546
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.033792281592267
$ws.Range("I16").Value = 1.030972391193527
$ws.Range("J16").Value = 1.032028494618229
$ws.Range("K16").Value = 1.032566522808608
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.037237680937575
$ws.Range("N16").Value = 1.03349409345591
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026258992327354
$ws.Range("D17").Value = 1.029553251984735
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.034444867276625
$ws.Range("I17").Value = 1.03111293358686
$ws.Range("J17").Value = 1.032506738099716
$ws.Range("K17").Value = 1.032941714378441
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.037816136235337
$ws.Range("N17").Value = 1.033973016097996
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026619859035813
$ws.Range("D18").Value = 1.029814864036049
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.03482521603226
$ws.Range("I18").Value = 1.031194625046711
$ws.Range("J18").Value = 1.0327853335672
$ws.Range("K18").Value = 1.033160209384191
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.038153195456741
$ws.Range("N18").Value = 1.034252007202998
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026742850121464
$ws.Range("D19").Value = 1.029904019425503
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.034954855777196
$ws.Range("I19").Value = 1.031222431567065
$ws.Range("J19").Value = 1.032880267374028
$ws.Range("K19").Value = 1.033234651768917
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.038268065980337
$ws.Range("N19").Value = 1.034347075826722
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026192587344926
$ws.Range("D20").Value = 1.029505107781585
$ws.Range("E20").Value = 0.9894336180355766
$ws.Range("F20").Value = 1.034374881460669
$ws.Range("I20").Value = 1.031097884175135
$ws.Range("J20").Value = 1.032455463984021
$ws.Range("K20").Value = 1.03290149591582
$ws.Range("L20").Value = 0.9929783193490043
$ws.Range("M20").Value = 1.03775410912628
$ws.Range("N20").Value = 1.033921669167176
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024399750707602
$ws.Range("D21").Value = 1.028204874925067
$ws.Range("E21").Value = 0.9882828385668255
$ws.Range("F21").Value = 1.032485848652709
$ws.Range("I21").Value = 1.030689616290198
$ws.Range("J21").Value = 1.031070181101726
$ws.Range("K21").Value = 1.031814267039254
$ws.Range("L21").Value = 0.9920501090198107
$ws.Range("M21").Value = 1.036079114334926
$ws.Range("N21").Value = 1.032534419024294
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023269479857593
$ws.Range("D22").Value = 1.027384768325802
$ws.Range("E22").Value = 0.9875604150241496
$ws.Range("F22").Value = 1.031295398943178
$ws.Range("I22").Value = 1.03043033415122
$ws.Range("J22").Value = 1.030195921039083
$ws.Range("K22").Value = 1.031127493347557
$ws.Range("L22").Value = 0.991467000034148
$ws.Range("M22").Value = 1.035022792138787
$ws.Range("N22").Value = 1.031658917412082
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023868950835941
$ws.Range("D23").Value = 1.027819771542369
$ws.Range("E23").Value = 0.9879432794636459
$ws.Range("F23").Value = 1.031926743254572
$ws.Range("I23").Value = 1.030568031018714
$ws.Range("J23").Value = 1.030659697030915
$ws.Range("K23").Value = 1.031491870139879
$ws.Range("L23").Value = 0.9917760702887607
$ws.Range("M23").Value = 1.035583074840069
$ws.Range("N23").Value = 1.032123352019018
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02622259390992
$ws.Range("D24").Value = 1.029526862945809
$ws.Range("E24").Value = 0.9894529299347241
$ws.Range("F24").Value = 1.034406505945065
$ws.Range("I24").Value = 1.031104685235741
$ws.Range("J24").Value = 1.032478633649453
$ws.Range("K24").Value = 1.032919669982897
$ws.Range("L24").Value = 0.9929938892766438
$ws.Range("M24").Value = 1.037782137569758
$ws.Range("N24").Value = 1.033944871736189
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028941066782685
$ws.Range("D25").Value = 1.031496818303435
$ws.Range("E25").Value = 0.9912096547607046
$ws.Range("F25").Value = 1.037272669910775
$ws.Range("I25").Value = 1.031716292813066
$ws.Range("J25").Value = 1.034575484359844
$ws.Range("K25").Value = 1.034562932959373
$ws.Range("L25").Value = 0.9944092447426411
$ws.Range("M25").Value = 1.040320578192477
$ws.Range("N25").Value = 1.036044700215102
